$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.788.27"
$ws.Range("E2").Value = "  -0.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.799.21"
$ws.Range("E3").Value = "  -1.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "702.56"

# Row 6 - Solana
Set-TextValue "D6" "169.77"
$ws.Range("E6").Value = "  -2.54%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.798.07"
$ws.Range("E7").Value = "  -1.85%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.83%  "

# Row 11 - Toncoin
Set-TextValue "D11" "7.62"
$ws.Range("E11").Value = "  +6.45%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.51%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -3.96%  "

# Row 14 - Avalanche
Set-TextValue "D14" "35.70"
$ws.Range("E14").Value = "  -2.25%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.438.07"
$ws.Range("E15").Value = "  -1.92%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.808.75"
$ws.Range("E16").Value = "  -1.52%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "70.763.82"
$ws.Range("E17").Value = "  -0.68%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  +0.32%  "

# Row 19 - was Chainlink, becomes Polkadot
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D19" "7.09"
$ws.Range("E19").Value = "  -2.25%  "

# Row 20 - was Polkadot, becomes Chainlink
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "17.32"
$ws.Range("E20").Value = "  -2.44%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "496.13"
$ws.Range("E21").Value = "  -0.91%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -4.77%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.722"
$ws.Range("E23").Value = "  -0.45%  "

# Row 24 - Litecoin
Set-TextValue "D24" "84.18"

# Row 25 - PEPE
$ws.Range("E25").Value = "  -4.47%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.948.94"
$ws.Range("E26").Value = "  -1.63%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.02"
$ws.Range("E27").Value = "  -2.21%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.26"
$ws.Range("E28").Value = "  -5.57%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  -7.22%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -5.98%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -0.67%  "

# Row 33 - NEARProtocol
Set-TextValue "D33" "7.31"
$ws.Range("E33").Value = "  -3.85%  "

# Row 34 - EthereumClassic
$ws.Range("E34").Value = "  -2.82%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.175"
$ws.Range("E35").Value = "  -3.31%  "

# Row 36 - was Binance-PegBSC-USD, becomes RenzoRestakedETH
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.765.68"
$ws.Range("E36").Value = "  -1.53%  "

# Row 37 - was RenzoRestakedETH, becomes Binance-PegBSC-USD
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D37" "0.998"
$ws.Range("E37").Value = "  -0.31%  "

# Row 38 - Aptos
Set-TextValue "D38" "9.04"
$ws.Range("E38").Value = "  -2.63%  "

# Row 39 - Hedera
$ws.Range("E39").Value = "  -3.69%  "

# Row 40 - Stacks
Set-TextValue "D40" "2.37"
$ws.Range("E40").Value = "  -1.65%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  -2.77%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  -1.81%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  -6.54%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.14%  "

# Row 46 - Monero
Set-TextValue "D46" "166.81"
$ws.Range("E46").Value = "  +1.90%  "

# Row 47 - FLOKI
$ws.Range("E47").Value = "  +0.41%  "

# Row 48 - OKB
Set-TextValue "D48" "48.91"
$ws.Range("E48").Value = "  -0.01%  "

# Row 49 - Bittensor
Set-TextValue "D49" "420.69"
$ws.Range("E49").Value = "  +0.36%  "

# Row 50 - Cosmos
Set-TextValue "D50" "8.57"
$ws.Range("E50").Value = "  -1.36%  "

# Row 51 - TheGraph
Set-TextValue "D51" "0.292"
$ws.Range("E51").Value = "  -3.89%  "
